$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (exact text, incl. leading/trailing spaces)
$updates = [ordered]@{
    "D2" = "43.054.18"
    "E2" = "  -4.70%  "
    "D3" = "2.225.28"
    "E3" = "  -5.59%  "
    "E4" = "  +0.09%  "
    "D5" = "318.25"
    "E5" = "  +2.81%  "
    "D6" = "98.86"
    "E6" = "  -8.46%  "
    "D7" = "0.581"
    "E7" = "  -7.26%  "
    "E8" = "  +0.03%  "
    "E9" = "  -7.99%  "
    "D10" = "36.53"
    "E10" = "  -11.06%  "
    "D11" = "54.41"
    "E11" = "  -1.87%  "
    "D12" = "0.0827"
    "E12" = "  -9.85%  "
    "D13" = "7.73"
    "E13" = "  -8.38%  "
    "E14" = "  -3.98%  "
    "D15" = "0.867"
    "E15" = "  -11.69%  "
    "D16" = "2.564.11"
    "E16" = "  -5.61%  "
    "E17" = "  -8.30%  "
    "D18" = "2.222.19"
    "E18" = "  -5.83%  "
    "D19" = "42.881.46"
    "E19" = "  -5.04%  "
    "D20" = "14.50"
    "E20" = "  +3.80%  "
    "D21" = "0.0₃0967"
    "E21" = "  -8.96%  "
    "D22" = "6.47"
    "E22" = "  -11.62%  "
    "D23" = "65.26"
    "E23" = "  -10.80%  "
    "D24" = "3.17"
    "E24" = "  -9.03%  "
    "D25" = "236.17"
    "E25" = "  -8.94%  "
    "E26" = "  -8.65%  "
    "E27" = "  +0.33%  "
    "D28" = "10.15"
    "E28" = "  -8.72%  "
    "D29" = "2.20"
    "E29" = "  -6.55%  "
    "E30" = "  -13.86%  "
    "D31" = "0.0884"
    "E31" = "  -8.37%  "
    "D32" = "20.52"
    "E32" = "  -7.90%  "
    "B33" = "InjectiveProtocol"
    "C33" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D33" = "33.79"
    "E33" = "  -10.50%  "
    "B34" = "Monero"
    "C34" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D34" = "156.57"
    "E34" = "  -7.30%  "
    "E35" = "  -5.58%  "
    "D36" = "3.32"
    "E36" = "  +11.81%  "
    "E37" = "  +17.64%  "
    "E38" = "  -5.96%  "
    "D39" = "4.47"
    "E39" = "  -6.90%  "
    "E40" = "  -11.28%  "
    "D41" = "3.67"
    "E41" = "  -6.26%  "
    "E42" = "  -8.66%  "
    "D43" = "1.869.49"
    "E43" = "  +11.39%  "
    "E44" = "  +0.07%  "
    "D45" = "12.09"
    "E45" = "  -5.94%  "
    "D46" = "88.00"
    "E46" = "  -11.47%  "
    "B47" = "Algorand"
    "C47" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D47" = "0.207"
    "E47" = "  -10.82%  "
    "B48" = "THORChain"
    "C48" = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
    "D48" = "5.50"
    "E48" = "  -0.08%  "
    "B49" = "ordi"
    "C49" = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
    "D49" = "78.57"
    "E49" = "  -4.16%  "
    "D50" = "60.60"
    "E50" = "  -12.86%  "
    "D51" = "8.61"
    "E51" = "  -5.70%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "1.00", "43.054.18", "0.581" keep their
    # exact literal representation instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
